$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.415.79"
$ws.Range("E2").Value = "  -4.49%  "

# Row 3
$ws.Range("D3").Value = "3.321.95"
$ws.Range("E3").Value = "  -5.05%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.88"
$ws.Range("E5").Value = "  -3.82%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.25"
$ws.Range("E6").Value = "  -4.15%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "3.324.54"
$ws.Range("E8").Value = "  -4.99%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  -0.84%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.37"
$ws.Range("E10").Value = "  -3.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.119"
$ws.Range("E11").Value = "  -4.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.377"
$ws.Range("E12").Value = "  -2.51%  "

# Row 13
$ws.Range("D13").Value = "3.881.16"
$ws.Range("E13").Value = "  -5.17%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.120"
$ws.Range("E14").Value = "  -0.38%  "

# Row 15
$ws.Range("D15").Value = "3.312.45"
$ws.Range("E15").Value = "  -5.29%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000169"
$ws.Range("E16").Value = "  -5.99%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.92"
$ws.Range("E17").Value = "  +0.29%  "

# Row 18
$ws.Range("D18").Value = "61.461.41"
$ws.Range("E18").Value = "  -4.33%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.53"
$ws.Range("E19").Value = "  +0.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.68"
$ws.Range("E20").Value = "  -1.32%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.07"
$ws.Range("E21").Value = "  -9.36%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "354.62"
$ws.Range("E22").Value = "  -8.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.556"
$ws.Range("E23").Value = "  -3.88%  "

# Row 24
$ws.Range("E24").Value = "  -0.05%  "

# Row 25
$ws.Range("D25").Value = "3.449.06"
$ws.Range("E25").Value = "  -5.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.47"
$ws.Range("E26").Value = "  -6.54%  "

# Row 27
$ws.Range("E27").Value = "  -6.87%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.17%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.19"
$ws.Range("E29").Value = "  -1.42%  "

# Row 30: 'InternetComputer(DFINITY)' -> 'Fetch.AI'
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.44"
$ws.Range("E30").Value = "  -4.04%  "

# Row 31: 'Fetch.AI' -> 'InternetComputer(DFINITY)'
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  -3.09%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.11"
$ws.Range("E32").Value = "  -5.99%  "

# Row 33: 'USDe' -> 'Kaspa'
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.150"
$ws.Range("E33").Value = "  -3.16%  "

# Row 34: 'Kaspa' -> 'USDe'
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").Value = "3.349.38"
$ws.Range("E35").Value = "  -5.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.64"
$ws.Range("E36").Value = "  -2.66%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.27"
$ws.Range("E37").Value = "  -2.60%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.83"
$ws.Range("E38").Value = "  -0.81%  "

# Row 39: 'Monero' -> 'ImmutableX'
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("E39").Value = "  -3.78%  "

# Row 40: 'ImmutableX' -> 'Monero'
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "160.32"
$ws.Range("E40").Value = "  -2.57%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0763"
$ws.Range("E41").Value = "  -2.61%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.21%  "

# Row 43: 'OKB' -> 'Filecoin'
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.39"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44: 'Filecoin' -> 'OKB'
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.00"
$ws.Range("E44").Value = "  -2.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.744"
$ws.Range("E45").Value = "  -7.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.12"
$ws.Range("E46").Value = "  -5.23%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.56"
$ws.Range("E47").Value = "  -5.00%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.40"
$ws.Range("E48").Value = "  -7.81%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.72"
$ws.Range("E49").Value = "  -0.92%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.868"
$ws.Range("E50").Value = "  -5.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.26"
$ws.Range("E51").Value = "  +1.71%  "
